$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 115, shifting existing rows 115:123 down to 116:124
$ws.Rows.Item(115).Insert()

# Populate the new row 115 with data
$ws.Cells.Item(115, 1).Value = 7
$ws.Cells.Item(115, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(115, 3).Value = "Ñuble"
$ws.Cells.Item(115, 4).Value = 44461
$ws.Cells.Item(115, 4).NumberFormat = $ws.Cells.Item(116, 4).NumberFormat
$ws.Cells.Item(115, 5).Value = 16
$ws.Cells.Item(115, 6).Value = 100112003
$ws.Cells.Item(115, 7).Value = "Ajo"
$ws.Cells.Item(115, 8).Value = "Chino"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 60
$ws.Cells.Item(115, 11).Value = 15500
$ws.Cells.Item(115, 12).Value = 16000
$ws.Cells.Item(115, 13).Value = 15750
$ws.Cells.Item(115, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(115, 15).Value = "China"
$ws.Cells.Item(115, 16).Value = 1575
$ws.Cells.Item(115, 17).Value = 10
$ws.Cells.Item(115, 18).Value = "Hortaliza"
